$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "ref" column for the 3 remaining student rows (group G18
# renumbered from STD000560.. to STD000130..)
$ws.Range("A2").Value = "STD000130-PROJ1-G18"
$ws.Range("A3").Value = "STD000131-PROJ1-G18"
$ws.Range("A4").Value = "STD000132-PROJ1-G18"

# Remove the trailing student rows (5-8) entirely, keeping their
# formatting but clearing ref/first_name/last_name/email/entrance/sex
$ws.Range("A5:F8").ClearContents()

# Move the active selection to A4 (was A8)
$ws.Range("A4").Select() | Out-Null
